$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: time of data update
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 13:05"

# Alemania (row 11) - updated case numbers
$ws.Range("B11").Value = 179153
$ws.Range("C11").Value = 132
$ws.Range("E11").Value = 11839
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 8314

# Rows 82/83: Bosnia y Herzegovina overtakes Tayikistan in ranking
$ws.Range("A82").Value = "Bosnia y Herzegovina"
$ws.Range("B82").Value = 2372
$ws.Range("C82").Value = 22
$ws.Range("D82").Value = 1614
$ws.Range("E82").Value = 617
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 141

$ws.Range("A83").Value = "Tayikistan"
$ws.Range("D83").Value = 1008
$ws.Range("E83").Value = 1298
$ws.Range("H83").Value = 44

# Row 126 (Malta) - updated case numbers
$ws.Range("B126").Value = 600
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 469

# Row 131 (Nepal) - updated case numbers
$ws.Range("B131").Value = 507
$ws.Range("C131").Value = 50
$ws.Range("D131").Value = 70
$ws.Range("E131").Value = 434

# Rows 137/138/139: Etiopia overtakes Estado de Palestina and Madagascar in ranking
$ws.Range("A137").Value = "Etiopia"
$ws.Range("B137").Value = 429
$ws.Range("C137").Value = 30
$ws.Range("D137").Value = 128
$ws.Range("E137").Value = 296
$ws.Range("H137").Value = 5

$ws.Range("A138").Value = "Estado de Palestina"
$ws.Range("B138").Value = 423
$ws.Range("D138").Value = 346
$ws.Range("E138").Value = 75

$ws.Range("A139").Value = "Madagascar"
$ws.Range("B139").Value = 405
$ws.Range("D139").Value = 131
$ws.Range("E139").Value = 272
$ws.Range("H139").Value = 2
